$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reorder the "Periodo Mora" values so that the shared-string table ends up
# with 1608, 1606, 1605, 1604 (in that order) and the associated "Valor Mora"
# amounts follow the period that was moved to the top (1608 -> 9193, the rest -> 27600).
$ws.Range("E16").Value = "1608"
$ws.Range("F16").Value = 9193

$ws.Range("E17").Value = "1606"
$ws.Range("F17").Value = 27600

$ws.Range("E18").Value = "1605"
$ws.Range("F18").Value = 27600

$ws.Range("E19").Value = "1604"
$ws.Range("F19").Value = 27600
